$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Host "Sheet name:" $ws.Name
Write-Host "H7:" $ws.Range("H7").Value
